$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.242.54"
$ws.Range("E2").Value = "  +5.48%  "

$ws.Range("D3").Value = "2.758.97"
$ws.Range("E3").Value = "  +3.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.56"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.53"
$ws.Range("E6").Value = "  +6.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +1.72%  "

$ws.Range("D9").Value = "2.756.04"
$ws.Range("E9").Value = "  +3.18%  "

$ws.Range("E10").Value = "  +2.02%  "

$ws.Range("E11").Value = "  +4.75%  "

$ws.Range("E12").Value = "  +2.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.160"
$ws.Range("E13").Value = "  +3.94%  "

$ws.Range("D14").Value = "3.219.32"
$ws.Range("E14").Value = "  +2.44%  "

$ws.Range("E15").Value = "  +2.26%  "

$ws.Range("D16").Value = "64.082.26"
$ws.Range("E16").Value = "  +5.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000153"
$ws.Range("E17").Value = "  +6.25%  "

$ws.Range("D18").Value = "2.752.37"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("E19").Value = "  +2.99%  "

$ws.Range("E20").Value = "  +2.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.11"
$ws.Range("E21").Value = "  +2.94%  "

$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.44"
$ws.Range("E25").Value = "  +3.77%  "

$ws.Range("E26").Value = "  +5.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  +4.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("E29").Value = "  +12.99%  "

$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("E31").Value = "  +4.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.26"
$ws.Range("E32").Value = "  +18.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.00"
$ws.Range("E33").Value = "  +3.33%  "

$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("E35").Value = "  +2.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.81"
$ws.Range("E36").Value = "  +7.85%  "

$ws.Range("E37").Value = "  +8.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +10.12%  "

$ws.Range("E39").Value = "  +16.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "344.81"
$ws.Range("E40").Value = "  +3.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  +5.45%  "

$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("E43").Value = "  +9.09%  "

$ws.Range("E44").Value = "  +6.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.77"
$ws.Range("E45").Value = "  +6.85%  "

$ws.Range("E46").Value = "  +5.00%  "

$ws.Range("E47").Value = "  +5.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.63"
$ws.Range("E48").Value = "  +2.80%  "

$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("E51").Value = "  -0.15%  "

